$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "Profilyyyyyyyyyyyy"
$ws.Range("B1").Select()
